$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Nguyễn Duy Cường" used to sit alone on row 22 (column D). Cường has now
# filled in his contact details, so his name moves down to row 23 and is
# joined by a phone number (E23) and an emailed hyperlink (F23).
$ws.Range("D22").Cut($ws.Range("D23"))

$ws.Range("E23").Value = 939025482

$ws.Range("F23").Value = "nhoxsazd@gmail.com"
$ws.Hyperlinks.Add($ws.Range("F23"), "mailto:nhoxsazd@gmail.com")

$ws.Range("F23").Select()
